$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.856.01"
$ws.Range("E2").Value = "  +1.89%  "

$ws.Range("D3").Value = "3.160.27"
$ws.Range("E3").Value = "  +3.87%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.38%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.159.75"
$ws.Range("E8").Value = "  +3.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.90%  "

$ws.Range("D15").Value = "3.678.54"
$ws.Range("E15").Value = "  +3.94%  "

$ws.Range("D16").Value = "64.946.37"
$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.76%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.147.97"
$ws.Range("E18").Value = "  +3.57%  "

$ws.Range("E19").Value = "  +1.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.725"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.51%  "

$ws.Range("E23").Value = "  +5.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.02%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.44%  "

$ws.Range("E31").Value = "  +13.92%  "

$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0893"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "471.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0419"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.49%  "

$ws.Range("D42").Value = "3.054.88"
$ws.Range("E42").Value = "  +1.78%  "

$ws.Range("E43").Value = "  +1.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.283"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.45%  "

$ws.Range("E45").Value = "  +8.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.17%  "

$ws.Range("D47").Value = "0.0₃0597"
$ws.Range("E47").Value = "  +17.04%  "

$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.114"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.17%  "
